# journal_de_bord.docx — add the "Jeudi 04 juin 2020 – 8eme jour" entry
# (commit: "fin du huitieme jour")
#
# Two changes:
#   1. The header shows a cached TIME field ("3 juin 2020") that must be
#      refreshed to "4 juin 2020" now that a new day was appended.
#   2. A page break plus the full 8th-day log (one Heading2 title, four
#      Heading3 time stamps and their Normal-style notes) is appended at
#      the very end of the document body, right before the sectPr.

$d = $word.ActiveDocument

# --- 1. Header date -------------------------------------------------------
$headerRange = $d.Sections(1).Headers(1).Range
$headerRange.Find.Execute(
    "3 juin 2020", $false, $false, $false, $false, $false,
    $true, 1, $false, "4 juin 2020", 2) | Out-Null

# --- 2. New day entry -------------------------------------------------------
# Create an empty paragraph at the end of the body to use as an insertion
# point, then fill it (and the paragraphs after it) via raw WordprocessingML
# so the run/paragraph formatting matches exactly.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null
$placeholder = $d.Paragraphs.Last

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newDayXml =
  "<w:p $w><w:r><w:br w:type=`"page`"/></w:r></w:p>" +
  "<w:p $w>" +
    "<w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" +
    "<w:r><w:lastRenderedPageBreak/><w:t>Jeudi 04 juin 2020 – 8</w:t></w:r>" +
    "<w:r><w:rPr><w:vertAlign w:val=`"superscript`"/></w:rPr><w:t>ème</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`"> jour</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:pPr><w:pStyle w:val=`"Heading3`"/></w:pPr>" +
    "<w:r><w:t>08h00</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:r><w:t>Pour commencer la journée, j’ai modifié la page d’accueil afin qu’elle n’affiche que les neuf films les plus (et les mieux) notés. La requête nécessaire à la récupération de ces films m’a pris beaucoup de temps à faire, j’avais du mal à la créer.</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:pPr><w:pStyle w:val=`"Heading3`"/></w:pPr>" +
    "<w:r><w:t>09h30</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:r><w:t>J’ai implémenté la modification de film à mon projet. J’ai eu des problèmes au niveau de la méthode, cette tâche a donc pris plus de temps que prévu.</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:pPr><w:pStyle w:val=`"Heading3`"/></w:pPr>" +
    "<w:r><w:t>11h30</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:r><w:t xml:space=`"preserve`">À ce moment-là, j’ai </w:t></w:r>" +
    "<w:r><w:t>documenté mes modifications et apporter des améliorations à ma documentation technique.</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:pPr><w:pStyle w:val=`"Heading3`"/></w:pPr>" +
    "<w:r><w:t>13h30</w:t></w:r>" +
  "</w:p>" +
  "<w:p $w>" +
    "<w:r><w:t>Une fois les grosses fonctionnalités du site programmées, j’ai voulu peaufiner en ajoutant des détails (messages d’erreurs, données sauvegardées en cas d’erreur, affichage des images lors des sélections, etc.)</w:t></w:r>" +
  "</w:p>"

$placeholder.Range.InsertXML($newDayXml) | Out-Null

Write-Output "journal_de_bord: 8e jour ajoute, en-tete mise a jour"
